$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.520.58"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "'2.429.03"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'558.59"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'160.26"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Value = "'0.164"
$ws.Range("E9").Value = "  +8.70%  "
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -5.17%  "
$ws.Range("D13").Value = "'68.399.08"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "'2.870.93"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("D16").Value = "'23.07"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "'2.426.19"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "'10.47"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "'335.09"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'66.70"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'3.68"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'2.553.44"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'0.0₃0819"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'425.05"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "'159.63"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").Value = "'19.04"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "'1.49"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").Value = "'1.08"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "'2.04"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'131.52"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'3.34"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "'0.0713"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'0.481"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("E51").Value = "  +0.09%  "
